$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New defined names -----------------------------------------------
$wb.Names.Add("PUE.NUM.legajos", $ws.Range("A3:A5"))
$wb.Names.Add("PUE.SWITCH.booleanos", $ws.Range("I9:I12"))

# --- Column A: "LEGAJOS" header + legajo numbers ----------------------
$ws.Range("A2").Value = "LEGAJOS"
$ws.Range("A3").Value = 25407
$ws.Range("A4").Value = 26549
$ws.Range("A5").Value = 23403

# --- Row 2: machine headers --------------------------------------------
$ws.Range("J2").Value = "Maquina 1"
$ws.Range("K2").Value = "Maquina 2"
$ws.Range("L2").Value = "Maquina 3"

# --- Row 3: Velocidad ----------------------------------------------------
$ws.Range("I3").Value = "Velocidad"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 200
$ws.Range("L3").Value = 700

# --- Row 4: Uso ------------------------------------------------------------
$ws.Range("I4").Value = "Uso"
$ws.Range("J4").Value = 0.9
$ws.Range("K4").Value = 0.3
$ws.Range("L4").Value = 0.5

# --- Row 5: Anios ------------------------------------------------------------
$ws.Range("I5").Value = "Anios"
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 2

# --- Column I, rows 9-12: boolean switches ---------------------------------
$ws.Range("I9").Value = $false
$ws.Range("I9").Style = "Normal"
$ws.Range("I10").Value = $false
$ws.Range("I10").Style = "Normal"
$ws.Range("I11").Value = $true
$ws.Range("I11").Style = "Normal"
$ws.Range("I12").Value = $false
$ws.Range("I12").Style = "Normal"

# --- Selection moves to L5 --------------------------------------------------
[void]$ws.Range("L5").Select()
